$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '268.95'
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = '18'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '22.77'
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = '18'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '6.337'
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = '18'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.06167'
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = '18'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '3.644'
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = '18'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '6.671'
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = '18'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.375'
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = '18'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.8291'
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = '18'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.01371'
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = '18'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1608'
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = '18'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08314'
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = '18'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03553'
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = '18'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.03219'
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = '18'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.09319'
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = '18'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.898'
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = '18'
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = '18'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.04749'
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = '18'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.006368'
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = '18'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.005654'
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = '18'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.001077'
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = '18'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0001503'
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = '18'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.723'
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = '18'
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = '18'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.3335'
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = '18'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1233'
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = '18'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0002709'
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = '18'
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = '18'
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = '18'
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = '18'
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = '18'
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = '18'
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = '18'
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = '18'
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = '18'
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = '18'
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = '18'
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = '18'
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = '18'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04713'
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = '18'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007001'
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = '18'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1156'
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = '18'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.003523'
$ws.Range("E43").Value = '42CEJICEJI'
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = '18'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.01176'
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = '18'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00006270'
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = '18'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0009918'
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = '18'
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = '18'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.7245'
$ws.Range("E48").Value = '47CoinbaseStockTokenCOINWorstin24h'
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = '18'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.002311'
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = '18'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00001403'
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = '18'
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = '18'
